$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column number format from the last existing data row (A13)
# down onto the new rows (A14:A36) so the new dates keep the same style (s="1")
# as the rest of column A, instead of Excel inventing a brand-new style.
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A14:A36").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# New rows of data: date serial number, temperature reading
# (one reading, for 2024-10-04, is missing -> stored as the text "None",
# reusing the existing shared string used elsewhere in the sheet).
$newRows = @(
    @(14, 45566, 20),
    @(15, 45567, 24),
    @(16, 45568, 22),
    @(17, 45569, "None"),
    @(18, 45570, 13),
    @(19, 45571, 50),
    @(20, 45572, 51),
    @(21, 45573, 45),
    @(22, 45574, 30),
    @(23, 45575, 27),
    @(24, 45576, 26),
    @(25, 45577, 25),
    @(26, 45578, 24),
    @(27, 45579, 60),
    @(28, 45580, 66),
    @(29, 45581, 64),
    @(30, 45582, 62),
    @(31, 45583, 40),
    @(32, 45584, 30),
    @(33, 45585, 33),
    @(34, 45586, 32),
    @(35, 45587, 33),
    @(36, 45588, 55)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $dateSerial = $row[1]
    $temp = $row[2]
    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $temp
}

# Match the author's final selection
$ws.Range("D35").Select() | Out-Null
